# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the 4 worker/period rows of the "estado de cuenta" table:
#  - VERONICA ISABEL HAMBURGER ESTRADA moves up to rows 16-18, one row per
#    mora period (1906, 1907, 1908), with her updated "Salario Basico".
#  - MARIA PAULA GUERRERO CHALELA moves down to row 19 (period 1909) with
#    her updated "Salario Basico".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: VERONICA ISABEL HAMBURGER ESTRADA - periodo 1906
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1235043865"
$ws.Range("D16").Value = "VERONICA ISABEL HAMBURGER ESTRADA"
$ws.Range("E16").Value = "1906"
$ws.Range("F16").Value = 36000
$ws.Range("G16").Value = 2050000

# Row 17: VERONICA ISABEL HAMBURGER ESTRADA - periodo 1907
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235043865"
$ws.Range("D17").Value = "VERONICA ISABEL HAMBURGER ESTRADA"
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 36000
$ws.Range("G17").Value = 2050000

# Row 18: VERONICA ISABEL HAMBURGER ESTRADA - periodo 1908
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235043865"
$ws.Range("D18").Value = "VERONICA ISABEL HAMBURGER ESTRADA"
$ws.Range("E18").Value = "1908"
$ws.Range("F18").Value = 36000
$ws.Range("G18").Value = 2050000

# Row 19: MARIA PAULA GUERRERO CHALELA - periodo 1909
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047426912"
$ws.Range("D19").Value = "MARIA PAULA GUERRERO CHALELA"
$ws.Range("E19").Value = "1909"
$ws.Range("F19").Value = 97520
$ws.Range("G19").Value = 2438000
